# Update the "Förändrad" (Changed) date column (C) for rows 2-14
# from 45185 (2023-09-16) to 45204 (2023-10-05).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 14; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45185) {
        $cell.Value2 = 45204
    }
}
